# Apply "Penalty Reward System" forecast shift edits
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Forecast Comparison" - shift Week_Start_Date (col B) forward
# one week and update MyForecast (col D) values for rows 2-17.
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

# Keep dates stored as plain text (inline/shared string), not as Excel
# date serials, matching the original workbook's formatting.
$wsForecast.Range("B2:B17").NumberFormat = "@"

$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$newForecast = @(18, 18, 18, 18, 18, 18, 18, 18, 19, 20, 20, 19, 19, 19, 19, 19)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 2).Value = $newDates[$i]
    $wsForecast.Cells.Item($row, 4).Value = $newForecast[$i]
}

# ---------------------------------------------------------------------
# Sheet 2: "Summary" - update summary metrics
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B13").NumberFormat = "@"
$wsSummary.Range("B15").NumberFormat = "@"

$wsSummary.Range("B2").Value = "2024-02-18 to 2025-01-05"
$wsSummary.Range("B8").Value = "561 units"
$wsSummary.Range("B9").Value = "300"
$wsSummary.Range("B10").Value = "146"
$wsSummary.Range("B11").Value = "73"
$wsSummary.Range("B12").Value = "20"
$wsSummary.Range("B13").Value = "2025-03-16"
$wsSummary.Range("B14").Value = "18"
$wsSummary.Range("B15").Value = "2025-01-12"
